$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70 (shifts existing rows 70-132 down to 71-133)
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new weekly data point
$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 45225
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100108
$ws.Cells.Item(70, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(70, 9).Value = 100108004
$ws.Cells.Item(70, 10).Value = "Papaya"
$ws.Cells.Item(70, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 100
$ws.Cells.Item(70, 14).Value = 24000
$ws.Cells.Item(70, 15).Value = 24000
$ws.Cells.Item(70, 16).Value = 24000
$ws.Cells.Item(70, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(70, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(70, 19).Value = 2400
$ws.Cells.Item(70, 20).Value = 10
